$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "61.428.13"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "2.984.16"
$ws.Range("E3").Value = "  -1.45%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "597.32"
$ws.Range("E5").Value = "  +1.83%  "
Set-TextValue $ws.Range("D6") "144.47"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D8").Value = "2.983.58"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("E15").Value = "  +3.36%  "
$ws.Range("D16").Value = "3.474.01"
$ws.Range("E16").Value = "  -1.42%  "
Set-TextValue $ws.Range("D17") "6.91"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "61.420.60"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "2.984.28"
Set-TextValue $ws.Range("D20") "444.47"
$ws.Range("E20").Value = "  -4.42%  "
Set-TextValue $ws.Range("D21") "13.90"
$ws.Range("E21").Value = "  -1.17%  "
Set-TextValue $ws.Range("D22") "0.681"
$ws.Range("E22").Value = "  -1.35%  "
Set-TextValue $ws.Range("D23") "7.31"
$ws.Range("E23").Value = "  -2.60%  "
Set-TextValue $ws.Range("D24") "80.96"
$ws.Range("E24").Value = "  -1.09%  "
Set-TextValue $ws.Range("D25") "10.74"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("E26").Value = "  -4.58%  "
Set-TextValue $ws.Range("D27") "12.03"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +2.30%  "
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.23%  "
Set-TextValue $ws.Range("D31") "7.26"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  -3.85%  "
Set-TextValue $ws.Range("D33") "27.13"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "0.0₃0811"
$ws.Range("E35").Value = "  -0.32%  "
Set-TextValue $ws.Range("D36") "1.02"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -1.03%  "
Set-TextValue $ws.Range("D38") "50.18"
$ws.Range("E38").Value = "  -0.77%  "
Set-TextValue $ws.Range("D39") "8.94"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("E40").Value = "  -6.14%  "
$ws.Range("E41").Value = "  +8.77%  "
Set-TextValue $ws.Range("D42") "2.86"
$ws.Range("E42").Value = "  -3.85%  "
Set-TextValue $ws.Range("D43") "387.83"
$ws.Range("E43").Value = "  -3.06%  "
Set-TextValue $ws.Range("D44") "39.58"
$ws.Range("E44").Value = "  +5.55%  "
Set-TextValue $ws.Range("D45") "0.268"
$ws.Range("E45").Value = "  -3.80%  "
Set-TextValue $ws.Range("D46") "0.0348"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("D47").Value = "2.683.55"
$ws.Range("E47").Value = "  -2.77%  "
Set-TextValue $ws.Range("D48") "131.64"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -2.38%  "
Set-TextValue $ws.Range("D51") "0.107"
$ws.Range("E51").Value = "  -2.86%  "
